# Update the cryptos worksheet figures (price / 1h-volume-change columns)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Some "Price" figures (column D) are plain digit strings like
    # "238.22" or "1.000" that Excel would otherwise auto-convert to
    # numbers when assigned through .Value. Force them to stay text,
    # matching the original inline-string cells, then drop the
    # number-format override again so no extra style is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "30.265.33"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "1.886.36"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
Set-TextValue "D5" "238.22"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("E6").Value = "  +0.19%  "

# Row 7
$ws.Range("E7").Value = "  -1.83%  "

# Row 8
Set-TextValue "D8" "0.2828"
$ws.Range("E8").Value = "  -0.53%  "

# Row 9
Set-TextValue "D9" "0.06575"
$ws.Range("E9").Value = "  -1.77%  "

# Row 10
Set-TextValue "D10" "19.73"
$ws.Range("E10").Value = "  +4.93%  "

# Row 11
Set-TextValue "D11" "0.07773"
$ws.Range("E11").Value = "  +1.23%  "

# Row 12
Set-TextValue "D12" "98.05"
$ws.Range("E12").Value = "  -3.11%  "

# Row 13
$ws.Range("D13").Value = "1.887.89"
$ws.Range("E13").Value = "  -1.79%  "

# Row 14
$ws.Range("E14").Value = "  -2.22%  "

# Row 15
Set-TextValue "D15" "0.6678"
$ws.Range("E15").Value = "  -0.39%  "

# Row 16
Set-TextValue "D16" "282.73"
$ws.Range("E16").Value = "  +10.27%  "

# Row 17
$ws.Range("D17").Value = "30.265.98"
$ws.Range("E17").Value = "  -0.89%  "

# Row 18
Set-TextValue "D18" "1.000"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("E19").Value = "  -0.48%  "

# Row 20
$ws.Range("D20").Value = "2.133.67"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21
Set-TextValue "D21" "5.363"
$ws.Range("E21").Value = "  -0.85%  "

# Row 22
Set-TextValue "D22" "0.000007302"
$ws.Range("E22").Value = "  -2.39%  "

# Row 23
Set-TextValue "D23" "1.000"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
Set-TextValue "D24" "6.172"
$ws.Range("E24").Value = "  -2.10%  "

# Row 25
Set-TextValue "D25" "9.354"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
Set-TextValue "D26" "165.39"
$ws.Range("E26").Value = "  -1.80%  "

# Row 27
Set-TextValue "D27" "19.12"
$ws.Range("E27").Value = "  +0.64%  "

# Row 28
Set-TextValue "D28" "1.993"
$ws.Range("E28").Value = "  -3.12%  "

# Row 29
Set-TextValue "D29" "1.376"
$ws.Range("E29").Value = "  -0.13%  "

# Row 30
Set-TextValue "D30" "0.09742"
$ws.Range("E30").Value = "  -3.30%  "

# Row 31
Set-TextValue "D31" "4.454"
$ws.Range("E31").Value = "  -5.57%  "

# Row 32
Set-TextValue "D32" "1.484"
$ws.Range("E32").Value = "  -1.90%  "

# Row 33
$ws.Range("E33").Value = "  -2.20%  "

# Row 34
Set-TextValue "D34" "0.04697"
$ws.Range("E34").Value = "  -0.64%  "

# Row 35
Set-TextValue "D35" "0.7072"
$ws.Range("E35").Value = "  -2.96%  "

# Row 36
$ws.Range("E36").Value = "  -1.57%  "

# Row 37
Set-TextValue "D37" "2.713"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
Set-TextValue "D38" "0.01866"
$ws.Range("E38").Value = "  -2.55%  "

# Row 39
Set-TextValue "D39" "6.664"
$ws.Range("E39").Value = "  +7.06%  "

# Row 40
Set-TextValue "D40" "2.525"
$ws.Range("E40").Value = "  -3.37%  "

# Row 41
Set-TextValue "D41" "72.08"
$ws.Range("E41").Value = "  -4.09%  "

# Row 42
Set-TextValue "D42" "0.8704"
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
Set-TextValue "D43" "1.969"
$ws.Range("E43").Value = "  +0.64%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D44" "1.001"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D45" "103.98"
$ws.Range("E45").Value = "  -1.32%  "

# Row 46
Set-TextValue "D46" "0.4198"
$ws.Range("E46").Value = "  -1.07%  "

# Row 47
Set-TextValue "D47" "985.65"
$ws.Range("E47").Value = "  -0.49%  "

# Row 48
Set-TextValue "D48" "7.214"
$ws.Range("E48").Value = "  -2.64%  "

# Row 49
Set-TextValue "D49" "9.247"
$ws.Range("E49").Value = "  +5.04%  "

# Row 50
Set-TextValue "D50" "0.1162"
$ws.Range("E50").Value = "  -3.04%  "

# Row 51
Set-TextValue "D51" "34.08"
$ws.Range("E51").Value = "  -2.17%  "
